# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps to reflect a fresh report run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (also shared by de-de!H2)
$wsOverview.Range("G2").Value = "2016-08-29 23:07:57"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-29 23:07:53"
$wsZhCn.Range("K2").Value = "2016-08-29 23:08:16"

# de-de sheet: Correspond Handback DateTime (H2 shares the same string as Overview!G2)
$wsDeDe.Range("H2").Value = "2016-08-29 23:07:57"
$wsDeDe.Range("K2").Value = "2016-08-29 23:08:23"
